$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The tutorial paragraphs describing the On-the-far-right-hand-side /
# Execute-Abort / Reset / clear button / speed dial / bytes / stop button /
# move-forward cluster were relocated later in the narrative and several
# of them received copyedits (typo fixes, punctuation, wording tweaks).
# Re-write the "Text Segment" cells (column C) of Sheet1 top to bottom so
# that the shared-string table compacts/re-appends in the same order as
# the target workbook.
$ws.Range("C7").Value = 'On the far right hand side there is a column of buttons.   Lets us begin at the beginning…or in this case, the top.  The {Execute} button will cause Rufus to begin to carry out all of your commands.  Once you push the {Execute button} you will notice that in changes to {Abort}.  You can push Abort at any time while Rufus is carrying out your commands to stop him from continuing to execute your list. Why would you need to do this you ask?  Well I would tell you if you quit interrupting.  The commands in your command window operate in a loop.  For all you non programmers out there this means that once you reach the end of your command list it will start back at the beginning and continue to execute commands from the top of your list.  '
$ws.Range("C8").Value = 'Below the Execute/Abort button you will find the {Reset} Button.  This button does just that, it resets Rufus to the start square.  It can even be used while Rufus is executing your commands.  '
$ws.Range("C9").Value = 'Now let’s say for a moment that you’ve completely blown it and you have no idea where you went wrong with your commands.  You have decided it’s just too much for your brain to handle and you want to start from scratch.  At this point you want to hit the {clear button}.  This will wipe your entire command window clean and allow you to start from scratch.  '
$ws.Range("C10").Value = 'The {help button} will take you to a help menu in game just in case you decided to skip this tutorial because you think you can figure it out on your own.  If you are not reading this you will have figure it out by yourself.'
$ws.Range("C11").Value = 'The {options button} will take you to the options menu.  This menu will allow you to do a number of things ranging from returning to the Main Menu to Skipping the remainder of this tutorial.  If you are not reading this you will have to figure it by yourself.  '
$ws.Range("C12").Value = 'Take note of the 600 MHZ with right and left arrows located above it.  This controls the speed and which Rufus will execute your commands.  Some of the levels will require you to play parts of them over again in order for you to know what to do next.  You will need to stop and add on commands once you have seen where Rufus ends up on the map.  You can speed up through parts you have already been through in the level or slow down in new parts to ensure you have put in the proper execution commands to solve the level.  '
$ws.Range("C13").Value = 'One last, yet very important, thing to bring to your attention is the small window above the map navigation tool.  You will notice that there are a number of bytes that are available for use to navigate the map and get to the end square.  Each command you issue Rufus takes up a certain number of you available bytes.  If you run out of bytes before you get Rufus to the end of the map you are working on you will have to go back and find a more efficient solution.   Subroutines can be very useful for this but let’s not get ahead of ourselves just yet.  For now know that there are a number of bytes available and the fewer bytes used, or more efficient your solution, the more points you will be awarded at the end of the level.  '
$ws.Range("C14").Value = 'Enough talk… well at least enough talk about all that other stuff.  Let’s get down to the buttons.  Before we talk about the buttons that make Rufus move we want to talk about the one button that does just the opposite, it makes Rufus stop in his tracks.  The {stop button} button can be placed into the command window at any point to bring Rufus to a stop in the execution even if there are commands located after the stop.  This command can be useful when you suddenly die because you put the wrong set of execution commands in and are now watching, with tears in your eyes, wondering where it all went wrong.  It can also help you break up longer levels into small parts.  You can place a stop command at the end of a small set of commands to ensure you navigated through that portion of the level correctly.  From there you can continue to add commands after the {stop button} and then remove it to see if you got the following section correct; and so on, and so forth. Remember with great power comes great responsibility so abuse the hell out of this and use it often.  It can help you during more difficult aspects of the game.  '
$ws.Range("C15").Value = 'You will not need to use the {stop button} during this tutorial process so you can get those ideas out of your head.  After small sections of the map we will guide you over each section and clear out your command window once you have completed it.  This way you will not need to run back through prior sections of the map and you can focus on the task in front of you.  During the regular game you will not have my powers of wisdom and you will have to figure it out on your own so soak it up while you can, Jonny (or if you''re name''s not Jonny, trace(yourName);).  '
$ws.Range("C16").Value = 'The first three buttons to focus on are the {move forward}, {turn left} and {turn right} buttons.  The {move forward} button moves Rufus forward one square.  The {turn left} button will rotate Rufus 90 degrees so he is now facing to the left.  Likewise the {turn right} button will turn Rufus to the right.  Now two wrongs don’t make a right but three rights will make you face left.  With this new found knowledge try moving Rufus through the S-curve in the first 5 tiles on the map. {End tutorial 1}'
$ws.Range("C17").Value = 'Congratulations, you made it through the first section of the map and proved that you are in fact a carbon based life form.  Let’s keep expanding our primitive minds with the next command.  The {move forward until unable} button is a button that can be very useful to cover long strait distances without having to push {move forward} many times over.  The {move forward until unable} button will have Rufus continue to travel in the direction he is facing at the time of the buttons execution until he dies or is no longer able to move forward because something is obstructing his path; such as a wall or door or edge of the map (yes the world is flat but no you will not sail off the edge of the world).  This tool is useful but beware of over use as its byte cost is high.  If you find yourself out of bytes you may want to go back and take a few of these out of your command list.  For now though, go nuts and enjoy yourself. {end move 02 }'
$ws.Range("C18").Value = 'Congratulations!!! Now Rufus can pass a sobriety test if he gets pulled over by the Norton Police.  In the Next section you will see a wall and a complete block that have a bunch of cracks in it.  Now Rufus might now seem like much but his Kung Fu is strong.  Rufus has the ability to take these weak walls and blocks and smash through them.  To get him to do this use the {punch} button while standing in front of them.  So go ahead…free your mind, punch through that wall.  Don’t think you can…know you can.  {end movie 03}'
$ws.Range("C19").Value = 'Congratulations!!!  Your sensei would be most proud of you.  Ahead of Rufus is a group of raised tiles.  The {climb} button can help you navigate this obstacle.  You will need to use this button once for each step.  It is important to note that Rufus can only climb up one level above his current height but he can drop off of any stair tile regardless of its height.  So climb up using {climb} button and move forward to jump off the back end.'
$ws.Range("C20").Value = 'Congratulations!!!  While that was no Everest don’t let anyone take that moment away from.  You too are  a unique and special snowflake.  In front of you there are three new obstacles; a gap, a duct and an electric tile.  To get past these obstacles you will need to use the {jump} button and the {crouch} button.  Both gaps and electric tiles can be jumped over by using the {jump} button.  You can even jump from one raised tile to another raised tile of the same height.  The duct, in between the two jump obstacles, must be ducked under for Rufus to continue.  For this you will use the {crouch} button.  So jump, limbo, and jump your way through to the next section.  {End Movie 05}'
$ws.Range("C21").Value = 'The next three tiles in front of you can all need to be activated by Rufus.  You can do this by use the {activate} button each time Rufus comes in contact with them. The first obstacle we see is a door that needs to be opened before Rufus can walk through it.  Doors can be on any side of a tile so it is important that you are facing the same direction as the door in order to activate and open it up.  '
$ws.Range("C22").Value = 'The next tile {reprogrammable square} is a Reprogrammable Square.  This tile operates in a similar way that these tutorial stop spots do in that it resets this to be Rufus’ new starting location should you die or need to go back and reset your current solution.  It also refunds any used bytes up to this point in your solution.  In order to access this tiles feature you will need to activate it while standing on them or facing them if they are on the edge of the tile. {top right reprogrammable square}  Reprogrammable squares can be activated as many times as you like within the level.  Making good use of them can help you increase your score by keeping your used bytes in the solution to a minimum.  '
$ws.Range("C23").Value = 'The {center switch tile} just after the Reprogrammable Square must also be activated.  This switch tile can be used for another number of different things in the game.  It has the ability to unblock paths by taking down walls or turning off electric tiles.  In addition it can take and inactive End Square {inactive end tile} and turn it back to active.  It even has the ability to make hidden portions of the map appear, like other switches, which you may need in order to get Rufus through the map.  The switch we see here is in the center of the tile which means no matter which direction Rufus is facing he can activate it.  Switches can also be on the edge of the tiles {upper left switch} and you must be facing them in order to activate them.  Take your new found power and go forth and concur.  {end movie 06}'
$ws.Range("C24").Value = 'Every time you want out….they keep pulling you back in!!!  At this point I believe you are ready to face the trials, which for this game consists of Subroutines.  As stated earlier you have a limited amount of bytes to use to navigate Rufus through each map.  You may notice parts of the map, like the one ahead, in which there is a repeating pattern to them.  In this case it is move forward three tiles and turn left or right.  Rather than put those commands in over and over again, which can be very byte expensive, we can put the commands to be repeated into the Sub1 and Sub2 section.   Why do this you ask with a puzzled look on your face?  Because I said so…okay there is more to it than that.  When you place the commands into the one of the two Sub’s it will cost you the original byte cost of each command that is placed.  But in the Main portion of the command screen when you want to execute the entire portion of code you placed into the Sub’s it will only cost you two bytes each time you use it.'
$ws.Range("C25").Value = 'For this section of code click on Sub1 on the main screen and put either {Move Forward} three times followed by {turn left} or {move forward till unable} followed by {turn left}.  Then click Sub2 and put the same thing except instead of using the {turn left} button use the {turn right} button.  Then in the Main command window use {Sub1} and {Sub2} as they match up with the map.  Give it a shot. {end movie 07}'
$ws.Range("C26").Value = 'Congratulations!!! No really this time I mean it.  Subroutines are not an easy concept to understand so give yourself a pat on the back. Just make sure you don’t choke on your gum because if you do we will lose our newfound respect for you.   Ahead of you are a group of {ice tile}; these are ice tiles.  Ice tiles are unique in that they are slippery.  Once you pick a direction on ice you only have to hit {move forward} one and Rufus will continue to slide until he hits a non-ice tile or hits an obstruction of some kind like a wall, door or edge of the map (again yes the world is flat but no you will not sail off the edge).  '
$ws.Range("C27").Value = 'The tile beyond this {teleport tile} is a teleport tile.  Hmm I wonder what this tile does.  Yes you guessed it, it teleports Rufus to another location of the map, you may now take your Fields Medal. {end movie 08}'
$ws.Range("C28").Value = 'Well you made!!! What a long strange trip it’s been.  All that’s left is to move forward on the End Square and claim this map in the name of whatever country you represent.  I hope you have enjoyed this tutorial and remember if you forget anything that has been discussed you can always click the {help button} to answer any question while you are in game.  Enjoy the game and thanks for playing!!!'

# Row 13 now holds slightly longer text ("One last, yet very important, ...")
# and grows a touch taller.
$ws.Rows.Item(13).RowHeight = 150

# Restore the saved scroll position / selection on the sheet.
$win = $excel.ActiveWindow
$win.ScrollRow = 25
$win.ScrollColumn = 1
$ws.Range("C28").Select() | Out-Null
